$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (shared strings: volume number + date range) ---
$ws.Range("A8").Value = "Volume 32   Number  37"
$ws.Range("C9").Value = "Report Covering the Week  9/8/2025  Through  9/14/2025"

# --- Crime complaint statistics table updates (rows 14-30) ---
$ws.Range("F14").NumberFormat = "@"
$ws.Range("F14").Value = "0"
$ws.Range("M14").Value = 150
$ws.Range("C15").NumberFormat = "#,##0"
$ws.Range("C15").Value = 1
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "***.*"
$ws.Range("F15").Value = 3
$ws.Range("H15").Value = 200
$ws.Range("I15").Value = 16
$ws.Range("K15").Value = 45.454545454545
$ws.Range("L15").Value = -30.434782608695
$ws.Range("M15").Value = -27.272727272727
$ws.Range("N15").Value = -76.119402985074
$ws.Range("C16").Value = 4
$ws.Range("D16").Value = 2
$ws.Range("E16").Value = 100
$ws.Range("F16").Value = 18
$ws.Range("G16").Value = 11
$ws.Range("H16").Value = 63.636363636363
$ws.Range("I16").Value = 165
$ws.Range("J16").Value = 161
$ws.Range("K16").Value = 2.484472049689
$ws.Range("L16").Value = 13.013698630137
$ws.Range("M16").Value = -17.085427135678
$ws.Range("N16").Value = -74.458204334365
$ws.Range("C17").Value = 8
$ws.Range("D17").Value = 9
$ws.Range("E17").Value = -11.111111111111
$ws.Range("G17").Value = 40
$ws.Range("H17").Value = -25
$ws.Range("I17").Value = 323
$ws.Range("J17").Value = 351
$ws.Range("K17").Value = -7.977207977207
$ws.Range("L17").Value = -3.582089552238
$ws.Range("M17").Value = 35.146443514644
$ws.Range("N17").Value = -54.442877291960
$ws.Range("C18").Value = 4
$ws.Range("D18").Value = 1
$ws.Range("E18").Value = 300
$ws.Range("F18").Value = 12
$ws.Range("G18").Value = 8
$ws.Range("H18").Value = 50
$ws.Range("I18").Value = 92
$ws.Range("J18").Value = 89
$ws.Range("K18").Value = 3.370786516853
$ws.Range("L18").Value = -32.352941176470
$ws.Range("M18").Value = 3.370786516853
$ws.Range("N18").Value = -86.685962373371
$ws.Range("C19").Value = 9
$ws.Range("D19").Value = 11
$ws.Range("E19").Value = -18.181818181818
$ws.Range("F19").Value = 30
$ws.Range("G19").Value = 28
$ws.Range("H19").Value = 7.142857142857
$ws.Range("I19").Value = 274
$ws.Range("J19").Value = 254
$ws.Range("K19").Value = 7.874015748031
$ws.Range("L19").Value = -13.836477987421
$ws.Range("M19").Value = 30.476190476190
$ws.Range("N19").Value = 1.481481481481
$ws.Range("D20").Value = 3
$ws.Range("F20").Value = 2
$ws.Range("G20").Value = 5
$ws.Range("H20").Value = -60
$ws.Range("J20").Value = 62
$ws.Range("K20").Value = -16.129032258064
$ws.Range("L20").Value = -46.391752577319
$ws.Range("M20").Value = 36.842105263157
$ws.Range("N20").Value = -71.111111111111
$ws.Range("C21").Value = 26
$ws.Range("D21").Value = 26
$ws.Range("E21").Value = 0
$ws.Range("F21").Value = 95
$ws.Range("G21").Value = 93
$ws.Range("H21").Value = 2.150537634408
$ws.Range("I21").Value = 932
$ws.Range("J21").Value = 936
$ws.Range("K21").Value = -0.427350427350
$ws.Range("L21").Value = -12.241054613936
$ws.Range("M21").Value = 16.354556803995
$ws.Range("N21").Value = -64.195159431425
$ws.Range("C23").Value = 1
$ws.Range("D23").Value = 2
$ws.Range("E23").Value = -50
$ws.Range("F23").Value = 11
$ws.Range("G23").Value = 17
$ws.Range("H23").Value = -35.294117647058
$ws.Range("I23").Value = 156
$ws.Range("J23").Value = 196
$ws.Range("K23").Value = -20.408163265306
$ws.Range("L23").Value = -9.826589595375
$ws.Range("M23").Value = 34.482758620689
$ws.Range("C24").Value = 19
$ws.Range("D24").Value = 18
$ws.Range("E24").Value = 5.555555555555
$ws.Range("F24").Value = 77
$ws.Range("G24").Value = 60
$ws.Range("H24").Value = 28.333333333333
$ws.Range("I24").Value = 685
$ws.Range("J24").Value = 586
$ws.Range("K24").Value = 16.894197952218
$ws.Range("L24").Value = 11.201298701298
$ws.Range("M24").Value = 49.563318777292
$ws.Range("C25").Value = 6
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "***.*"
$ws.Range("F25").Value = 17
$ws.Range("G25").Value = 5
$ws.Range("H25").Value = 240
$ws.Range("I25").Value = 143
$ws.Range("K25").Value = 33.644859813084
$ws.Range("L25").Value = 15.322580645161
$ws.Range("C26").Value = 11
$ws.Range("D26").Value = 21
$ws.Range("E26").Value = -47.619047619047
$ws.Range("F26").Value = 48
$ws.Range("G26").Value = 56
$ws.Range("H26").Value = -14.285714285714
$ws.Range("I26").Value = 477
$ws.Range("J26").Value = 513
$ws.Range("K26").Value = -7.017543859649
$ws.Range("L26").Value = 13.571428571428
$ws.Range("M26").Value = -22.690437601296
$ws.Range("C27").NumberFormat = "#,##0"
$ws.Range("C27").Value = 1
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "***.*"
$ws.Range("F27").Value = 3
$ws.Range("H27").Value = 50
$ws.Range("I27").Value = 19
$ws.Range("K27").Value = 5.555555555555
$ws.Range("L27").Value = -48.648648648648
$ws.Range("C28").NumberFormat = "#,##0"
$ws.Range("C28").Value = 3
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "***.*"
$ws.Range("F28").Value = 7
$ws.Range("H28").Value = 0
$ws.Range("I28").Value = 52
$ws.Range("K28").Value = 40.540540540540
$ws.Range("L28").Value = 36.842105263157
$ws.Range("M29").Value = -58.064516129032
$ws.Range("N29").Value = -85.393258426966
$ws.Range("M30").Value = -61.538461538461
$ws.Range("N30").Value = -87.804878048780
